# Rename the Pearson / BTEC logo pictures that live in this document's
# headers and footers.
#
#   word/footer1.xml (the "first page" footer)  : image1.png -> image2.png
#   word/footer2.xml (the "default" footer)      : image1.png -> image2.png
#   word/header1.xml (the "first page" header)   : image2.jpg -> image1.jpg
#   word/header2.xml (the "default" header)      : image2.jpg -> image1.jpg
#
# The rename touches the InlineShape's Name (OOXML wp:docPr/@name, mirrored
# onto pic:cNvPr/@name) and nothing else about the picture (size, position,
# alt text, etc. are left untouched).

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Headers ---------------------------------------------------------
# wdHeaderFooterPrimary = 1 (the default header used on most pages)
$null = $sec.Headers.Item(1).Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.jpg"

# wdHeaderFooterFirstPage = 2 (the header used on the first page)
$null = $sec.Headers.Item(2).Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.jpg"

# --- Footers ---------------------------------------------------------
# wdHeaderFooterPrimary = 1 (the default footer used on most pages)
$null = $sec.Footers.Item(1).Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# wdHeaderFooterFirstPage = 2 (the footer used on the first page)
$null = $sec.Footers.Item(2).Range.InlineShapes.Item(1).Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

Write-Output "Renamed header/footer logo pictures."
